$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New record row 9 (mirrors the existing row layout, columns per header row 1)
$ws.Range("A9").Value = 131302060
$ws.Range("B9").Value = 96605
$ws.Range("D9").Value = "LC"
$ws.Range("E9").Value = 2180
$ws.Range("F9").Value = "Blåmossa"
$ws.Range("G9").Value = "Leucobryum glaucum"
$ws.Range("H9").Value = "(Hedw.) Ångstr."

$ws.Range("I9").NumberFormat = "@"
$ws.Range("I9").Value = "5"

$ws.Range("P9").Value = "Hjärtaboda, Hjärtaboda, Sk"
$ws.Range("Q9").Value = 448906
$ws.Range("R9").Value = 6230159
$ws.Range("S9").Value = 20
$ws.Range("T9").Value = "Skåne"
$ws.Range("U9").Value = "Östra Göinge"
$ws.Range("V9").Value = "Skåne"
$ws.Range("W9").Value = "Hjärsås"

$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2026-02-25"

$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2026-02-25"

$ws.Range("AC9").Value = "Ett ex 60 cm bred, 20-30 cm hög"

$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false

$ws.Range("AW9").Value = "Martin Kornhall"
$ws.Range("AX9").Value = "Martin Kornhall"
